$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: year label in column A, value in column B
# Old data occupied rows 2-22 (years 2000-2020). New data occupies rows 2-14
# (years 2010-2022). Clear the old range first, then write fresh values.

$oldLastRow = 22
$newLastRow = 14

# Clear out the old data range beyond what we will rewrite (rows 2 through oldLastRow)
# Use ClearContents (not Clear) so that cell formatting/styles are preserved
# for the rows that remain (2..newLastRow).
$ws.Range("A2:B$oldLastRow").ClearContents() | Out-Null

# Remove the now-unused trailing rows entirely so the sheet dimension shrinks
# back down to A1:B$newLastRow, matching the smaller dataset.
if ($oldLastRow -gt $newLastRow) {
    $ws.Range("A$($newLastRow + 1):B$oldLastRow").EntireRow.Delete() | Out-Null
}

$years = @(2010, 2011, 2012, 2013, 2014, 2015, 2016, 2017, 2018, 2019, 2020, 2021, 2022)
$values = @(
    3906.6,
    4763.5589397997,
    6437.0682977038,
    7469.1253674561,
    8577.178987232601,
    9835.789637071901,
    11406.9815696409,
    13424.2244691549,
    17697.4212963104,
    22398.3881623544,
    28251.5091667389,
    37294.30297,
    47791
)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "$($years[$i])年"
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
